# UPDATE technology portfolios for Norway
# Update p_wi_q_waste (B2) and p_wi_c_waste (C2) across all year sheets,
# and refresh the base-year (2025) p_wi_c_inv (H2) hard-coded value.
# The dependent sheets (2030-2050) keep their existing formulas in H2,
# which will recalculate automatically from the new 2025 value.

$wb = $excel.ActiveWorkbook

$sheetNames = @("2025", "2030", "2035", "2040", "2045", "2050")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)
    $ws.Range("B2").Value = 0.3606149659239804
    $ws.Range("C2").Value = 50
}

# Base year (2025) sheet H2 is a hard-coded value (not a formula).
$ws2025 = $wb.Worksheets.Item("2025")
$ws2025.Range("H2").Value = 9268394

# Column H on the 2025 sheet no longer needs to be best-fit-widened;
# reset it back toward the workbook's default column width (the closest
# value reachable through the ColumnWidth property's internal rounding).
$ws2025.Columns.Item(8).ColumnWidth = 7.83
